$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Septiembre de 2020 a las 14:13"

# --- Update case-count numbers (columns B-H) for affected rows ---
$ws.Range("B4").Value = 6390176
$ws.Range("C4").Value = 1119
$ws.Range("D4").Value = 3636283
$ws.Range("E4").Value = 2561747
$ws.Range("G4").Value = 35
$ws.Range("H4").Value = 192146

$ws.Range("B6").Value = 4034339
$ws.Range("C6").Value = 14100
$ws.Range("D6").Value = 3112669
$ws.Range("E6").Value = 851921
$ws.Range("G6").Value = 114
$ws.Range("H6").Value = 69749

$ws.Range("B15").Value = 384666
$ws.Range("C15").Value = 1894
$ws.Range("D15").Value = 332131
$ws.Range("E15").Value = 30381
$ws.Range("G15").Value = 110
$ws.Range("H15").Value = 22154

$ws.Range("B38").Value = 88963
$ws.Range("C38").Value = 720
$ws.Range("D38").Value = 79903
$ws.Range("E38").Value = 8520
$ws.Range("G38").Value = 3
$ws.Range("H38").Value = 540

$ws.Range("B42").Value = 77040
$ws.Range("D42").Value = 65029
$ws.Range("E42").Value = 9186
$ws.Range("H42").Value = 2825

$ws.Range("B43").Value = 73471
$ws.Range("C43").Value = 705
$ws.Range("D43").Value = 63652
$ws.Range("E43").Value = 9431
$ws.Range("G43").Value = 1
$ws.Range("H43").Value = 388

$ws.Range("B44").Value = 73208
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 6237

$ws.Range("B45").Value = 72485
$ws.Range("D45").Value = 71510
$ws.Range("E45").Value = 274
$ws.Range("H45").Value = 701

$ws.Range("B46").Value = 70387
$ws.Range("C46").Value = 567
$ws.Range("D46").Value = 52346
$ws.Range("E46").Value = 15928
$ws.Range("G46").Value = 13
$ws.Range("H46").Value = 2113

$ws.Range("B47").Value = 70268
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 60417
$ws.Range("E47").Value = 8521
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 1330

$ws.Range("B48").Value = 68605
$ws.Range("D48").Value = 52483
$ws.Range("E48").Value = 14830
$ws.Range("H48").Value = 1292

$ws.Range("B49").Value = 63798
$ws.Range("C49").Value = 640
$ws.Range("D49").Value = 12347
$ws.Range("E49").Value = 49467
$ws.Range("G49").Value = 30
$ws.Range("H49").Value = 1984

$ws.Range("B50").Value = 59457
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 42576
$ws.Range("E50").Value = 15048
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 1833

$ws.Range("B51").Value = 56982
$ws.Range("C51").Value = 34
$ws.Range("D51").Value = 56174
$ws.Range("E51").Value = 781
$ws.Range("H51").Value = 27

$ws.Range("B52").Value = 56516
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 20612
$ws.Range("E52").Value = 35024
$ws.Range("H52").Value = 880

$ws.Range("B53").Value = 54743
$ws.Range("D53").Value = 42816
$ws.Range("E53").Value = 10876
$ws.Range("H53").Value = 1051

$ws.Range("B54").Value = 54095
$ws.Range("D54").Value = 50323
$ws.Range("E54").Value = 3576
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 196

$ws.Range("B55").Value = 50973
$ws.Range("D55").Value = 41249
$ws.Range("E55").Value = 9312
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 412

$ws.Range("B56").Value = 45773
$ws.Range("D56").Value = 32259
$ws.Range("E56").Value = 11975
$ws.Range("H56").Value = 1539

$ws.Range("B57").Value = 45680
$ws.Range("D57").Value = 18053
$ws.Range("E57").Value = 27158
$ws.Range("H57").Value = 469

$ws.Range("B58").Value = 45277
$ws.Range("C58").Value = 1041
$ws.Range("D58").Value = 27127
$ws.Range("E58").Value = 17870
$ws.Range("G58").Value = 9
$ws.Range("H58").Value = 280

$ws.Range("B59").Value = 44777
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 43693
$ws.Range("E59").Value = 801
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 283

$ws.Range("B60").Value = 44649
$ws.Range("C60").Value = 188
$ws.Range("D60").Value = 39823
$ws.Range("E60").Value = 3931
$ws.Range("G60").Value = 4
$ws.Range("H60").Value = 895

$ws.Range("B61").Value = 44293
$ws.Range("C61").Value = 94
$ws.Range("D61").Value = 39599
$ws.Range("E61").Value = 3634
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 1060

$ws.Range("B62").Value = 43957
$ws.Range("C62").Value = 425
$ws.Range("D62").Value = 36500
$ws.Range("E62").Value = 5444
$ws.Range("H62").Value = 2013

$ws.Range("B63").Value = 43075
$ws.Range("C63").Value = 77
$ws.Range("D63").Value = 40432
$ws.Range("E63").Value = 2304
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 339

$ws.Range("B64").Value = 38906
$ws.Range("D64").Value = 27017
$ws.Range("E64").Value = 10842
$ws.Range("H64").Value = 1047

$ws.Range("B65").Value = 38324
$ws.Range("C65").Value = 20
$ws.Range("D65").Value = 30082
$ws.Range("E65").Value = 6833
$ws.Range("H65").Value = 1409

$ws.Range("B66").Value = 37031
$ws.Range("C66").Value = 0
$ws.Range("D66").Value = 34419
$ws.Range("E66").Value = 2069
$ws.Range("H66").Value = 543

$ws.Range("B67").Value = 34884
$ws.Range("D67").Value = 21059
$ws.Range("E67").Value = 13236
$ws.Range("H67").Value = 589

$ws.Range("B68").Value = 31772
$ws.Range("D68").Value = 30387
$ws.Range("E68").Value = 664
$ws.Range("H68").Value = 721

$ws.Range("B69").Value = 29303
$ws.Range("D69").Value = 23364
$ws.Range("E69").Value = 4162
$ws.Range("H69").Value = 1777

$ws.Range("B70").Value = 29087
$ws.Range("C70").Value = 358
$ws.Range("D70").Value = 24828
$ws.Range("E70").Value = 3524
$ws.Range("H70").Value = 735

$ws.Range("B71").Value = 27249
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 19027
$ws.Range("E71").Value = 7793
$ws.Range("H71").Value = 429

$ws.Range("B72").Value = 26207
$ws.Range("C72").Value = 71
$ws.Range("D72").Value = 22330
$ws.Range("E72").Value = 3129
$ws.Range("G72").Value = 11
$ws.Range("H72").Value = 748

$ws.Range("B73").Value = 26207
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 5835

$ws.Range("B76").Value = 21439
$ws.Range("C76").Value = 297
$ws.Range("D76").Value = 14636
$ws.Range("E76").Value = 6152
$ws.Range("G76").Value = 12
$ws.Range("H76").Value = 651

$ws.Range("B85").Value = 15269
$ws.Range("C85").Value = 82
$ws.Range("D85").Value = 14113
$ws.Range("E85").Value = 957
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 199

$ws.Range("B93").Value = 10178
$ws.Range("C93").Value = 29
$ws.Range("D93").Value = 9420
$ws.Range("E93").Value = 499

$ws.Range("D153").Value = 1237
$ws.Range("E153").Value = 244

# --- Re-rank countries row 42-73: shift names so Suecia moves to the bottom of the block ---
$ws.Range("A42").Value = "Guatemala"
$ws.Range("A43").Value = "Emiratos Arabes Unidos"
$ws.Range("A44").Value = "Paises Bajos"
$ws.Range("A45").Value = "Bielorrusia"
$ws.Range("A46").Value = "Polonia"
$ws.Range("A47").Value = "Japon"
$ws.Range("A48").Value = "Marruecos"
$ws.Range("A49").Value = "Honduras"
$ws.Range("A50").Value = "Portugal"
$ws.Range("A51").Value = "Singapur"
$ws.Range("A52").Value = "Etiopia"
$ws.Range("A53").Value = "Nigeria"
$ws.Range("A54").Value = "Barein"
$ws.Range("A55").Value = "Venezuela"
$ws.Range("A56").Value = "Argelia"
$ws.Range("A57").Value = "Costa Rica"
$ws.Range("A58").Value = "Nepal"
$ws.Range("A59").Value = "Ghana"
$ws.Range("A60").Value = "Armenia"
$ws.Range("A61").Value = "Kirguistan"
$ws.Range("A62").Value = "Suiza"
$ws.Range("A63").Value = "Uzbekistan"
$ws.Range("A64").Value = "Moldavia"
$ws.Range("A65").Value = "Afganistan"
$ws.Range("A66").Value = "Azerbaiyan"
$ws.Range("A67").Value = "Kenia"
$ws.Range("A68").Value = "Serbia"
$ws.Range("A69").Value = "Irlanda"
$ws.Range("A70").Value = "Austria"
$ws.Range("A71").Value = "Chequia"
$ws.Range("A72").Value = "Australia"
$ws.Range("A73").Value = "Suecia"
